$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 354
$ws.Range("I33").Value = 389.125
$ws.Range("K33").Value = 389.125
$ws.Range("M33").Value = -160.125
$ws.Range("H92").Value = 838.5185
$ws.Range("I92").Value = 891.2273
$ws.Range("J92").Value = 606.6
$ws.Range("K92").Value = 891.2273
$ws.Range("L92").Value = 606.6
$ws.Range("M92").Value = 356.7727
$ws.Range("N92").Value = -3102.6
$ws.Range("H138").Value = 3501.513
$ws.Range("I138").Value = 1936.75
$ws.Range("K138").Value = 5810.25
$ws.Range("M138").Value = -670.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 50000
$ws.Range("I44").Value = 50000
$ws.Range("K44").Value = 50000
$ws.Range("M44").Value = -49512
$ws.Range("H55").Value = 58053
$ws.Range("J55").Value = 58053
$ws.Range("L55").Value = 58053
$ws.Range("N55").Value = -58683
$ws.Range("H61").Value = 2724
$ws.Range("I61").Value = 2408.3
$ws.Range("K61").Value = 2408.3
$ws.Range("M61").Value = -2196.3
$ws.Range("H102").Value = 1542.5333
$ws.Range("I102").Value = 1542.5333
$ws.Range("K102").Value = 1542.5333
$ws.Range("M102").Value = 79.46669999999995
$ws.Range("H122").Value = 1837
$ws.Range("I122").Value = 1319.3334
$ws.Range("K122").Value = 3958.0002
$ws.Range("M122").Value = -1508.0002
$ws.Range("H136").Value = 2724
$ws.Range("I136").Value = 2408.3
$ws.Range("K136").Value = 7224.900000000001
$ws.Range("M136").Value = -4674.900000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H86").Value = 3007.1428
$ws.Range("I86").Value = 1818.6364
$ws.Range("J86").Value = 4314.5
$ws.Range("K86").Value = 1818.6364
$ws.Range("L86").Value = 4314.5
$ws.Range("M86").Value = -695.6364000000001
$ws.Range("N86").Value = -6560.5
$ws.Range("H89").Value = 3007.1428
$ws.Range("I89").Value = 1818.6364
$ws.Range("J89").Value = 4314.5
$ws.Range("K89").Value = 9093.182000000001
$ws.Range("L89").Value = 21572.5
$ws.Range("M89").Value = -3477.182000000001
$ws.Range("N89").Value = -32804.5
$ws.Range("H99").Value = 12089.889
$ws.Range("I99").Value = 13501.125
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 13501.125
$ws.Range("L99").Value = 800
$ws.Range("M99").Value = -12003.125
$ws.Range("N99").Value = -3796
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 39935
$ws.Range("J50").Value = 39935
$ws.Range("L50").Value = 39935
$ws.Range("N50").Value = -41185
$ws.Range("H60").Value = 28842
$ws.Range("J60").Value = 47681
$ws.Range("L60").Value = 47681
$ws.Range("N60").Value = -48703
$ws.Range("H141").Value = 236819
$ws.Range("J141").Value = 236819
$ws.Range("L141").Value = 236819
$ws.Range("N141").Value = -247179
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.608696
$ws.Range("I2").Value = 37.545456
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 225.272736
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -112.272736
$ws.Range("N2").Value = -1426
$ws.Range("H11").Value = 138.66667
$ws.Range("J11").Value = 97
$ws.Range("L11").Value = 291
$ws.Range("N11").Value = -571
$ws.Range("H33").Value = 43.9
$ws.Range("I33").Value = 52.625
$ws.Range("J33").Value = 9
$ws.Range("K33").Value = 315.75
$ws.Range("L33").Value = 54
$ws.Range("M33").Value = -32.75
$ws.Range("N33").Value = -620
$ws.Range("H55").Value = 10309.375
$ws.Range("J55").Value = 14395
$ws.Range("L55").Value = 43185
$ws.Range("N55").Value = -43539
$ws.Range("H68").Value = 1749
$ws.Range("J68").Value = 1749
$ws.Range("L68").Value = 5247
$ws.Range("N68").Value = -6869
$ws.Range("H71").Value = 1749
$ws.Range("J71").Value = 1749
$ws.Range("L71").Value = 15741
$ws.Range("N71").Value = -23853
$ws.Range("H86").Value = 499.33334
$ws.Range("J86").Value = 699
$ws.Range("L86").Value = 2097
$ws.Range("N86").Value = -4469
$ws.Range("H89").Value = 499.33334
$ws.Range("J89").Value = 699
$ws.Range("L89").Value = 6291
$ws.Range("N89").Value = -18147
$ws.Range("H98").Value = 10000
$ws.Range("I98").Value = 10000
$ws.Range("K98").Value = 30000
$ws.Range("M98").Value = -28502
$ws.Range("H103").Value = 2594
$ws.Range("J103").Value = 3218.75
$ws.Range("L103").Value = 9656.25
$ws.Range("N103").Value = -11414.25
$ws.Range("H107").Value = 449.5
$ws.Range("J107").Value = 699.5
$ws.Range("L107").Value = 2098.5
$ws.Range("N107").Value = -5938.5
$ws.Range("H109").Value = 4666.619
$ws.Range("I109").Value = 2000
$ws.Range("J109").Value = 4947.316
$ws.Range("K109").Value = 6000
$ws.Range("L109").Value = 14841.948
$ws.Range("M109").Value = -4960
$ws.Range("N109").Value = -16921.948
$ws.Range("H131").Value = 3288.087
$ws.Range("I131").Value = 2499
$ws.Range("J131").Value = 3323.9546
$ws.Range("K131").Value = 7497
$ws.Range("L131").Value = 9971.863799999999
$ws.Range("M131").Value = -2457
$ws.Range("N131").Value = -20051.8638
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2016.6666
$ws.Range("I122").Value = 2076.889
$ws.Range("J122").Value = 1836
$ws.Range("K122").Value = 6230.667
$ws.Range("L122").Value = 5508
$ws.Range("M122").Value = -3780.667
$ws.Range("N122").Value = -10408
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3464.6667
$ws.Range("I122").Value = 3694.5
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 11083.5
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -8633.5
$ws.Range("N122").Value = -13915
$ws.Range("H136").Value = 2196.3333
$ws.Range("J136").Value = 1777
$ws.Range("L136").Value = 5331
$ws.Range("N136").Value = -10431
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5985.5
$ws.Range("I136").Value = 5515.4287
$ws.Range("J136").Value = 7959.8
$ws.Range("K136").Value = 16546.2861
$ws.Range("L136").Value = 23879.4
$ws.Range("M136").Value = -13996.2861
$ws.Range("N136").Value = -28979.4
